# Apply the "prevention and rehab 2.0.0" update to the workbook.
$wb = $excel.ActiveWorkbook

# --- Rename the "Include from Tempcodes" sheet -> "Include from CareSocialCodes" ---
$wsInclude = $wb.Worksheets.Item("Include from Tempcodes")
$wsInclude.Name = "Include from CareSocialCodes"

# --- Update values on the Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 2.0.0
$wsMeta.Range("B3").Value = "2.0.0"

# Date: 2023-09-28T10:05:19+02:00 -> 2024-06-14T10:48:54+02:00
$wsMeta.Range("B8").Value = "2024-06-14T10:48:54+02:00"

# Contact: "No display for ContactDetail" -> "Kommunernes Landsforening (http://kl.dk)"
$wsMeta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- Update the System URI on the renamed include sheet ---
# http://fhir.kl.dk/rehab/CodeSystem/Tempcodes -> http://fhir.kl.dk/term/CodeSystem/CareSocialCodes
$wsInclude.Range("B8").Value = "http://fhir.kl.dk/term/CodeSystem/CareSocialCodes"
